# Generate Report for Archive
#
# Updates the localization status from "Ready for handoff" to
# "In Translation" wherever it appears (Overview!E2:F2, zh-cn!C2,
# de-de!C2), and narrows the now-shorter "Status" columns to match the
# new content width (Overview cols E/F, zh-cn col C, de-de col C).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status text -----------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# --- Resize the status columns to fit the new (shorter) text ----------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5
